$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 437.84
$ws.Range("I15").Value = 437.84
$ws.Range("K15").Value = 1313.52
$ws.Range("M15").Value = -1144.52
$ws.Range("H17").Value = 1116.7164
$ws.Range("I17").Value = 778.5714
$ws.Range("J17").Value = 1156.1666
$ws.Range("K17").Value = 2335.7142
$ws.Range("L17").Value = 3468.4998
$ws.Range("M17").Value = -2167.7142
$ws.Range("N17").Value = -3804.4998
$ws.Range("H19").Value = 2334.7693
$ws.Range("I19").Value = 3130.6667
$ws.Range("K19").Value = 3130.6667
$ws.Range("M19").Value = -2955.6667
$ws.Range("H98").Value = 1297.25
$ws.Range("I98").Value = 1182.1428
$ws.Range("J98").Value = 2103
$ws.Range("K98").Value = 1182.1428
$ws.Range("L98").Value = 2103
$ws.Range("M98").Value = 315.8571999999999
$ws.Range("N98").Value = -5099
$ws.Range("H122").Value = 1297.25
$ws.Range("I122").Value = 1182.1428
$ws.Range("J122").Value = 2103
$ws.Range("K122").Value = 3546.4284
$ws.Range("L122").Value = 6309
$ws.Range("M122").Value = -1096.4284
$ws.Range("N122").Value = -11209
$ws.Range("H127").Value = 401.91666
$ws.Range("I127").Value = 401.91666
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 1205.74998
$ws.Range("L127").Value = 0
$ws.Range("M127").Value = 3754.25002
$ws.Range("N127").ClearContents()
$ws.Range("H132").Value = 4199.8936
$ws.Range("I132").Value = 1938.0883
$ws.Range("J132").Value = 10115.385
$ws.Range("K132").Value = 5814.2649
$ws.Range("L132").Value = 30346.155
$ws.Range("M132").Value = -3284.2649
$ws.Range("N132").Value = -35406.155

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2792.2952
$ws.Range("I132").Value = 2508.465
$ws.Range("J132").Value = 3470.3333
$ws.Range("K132").Value = 7525.395
$ws.Range("L132").Value = 10410.9999
$ws.Range("M132").Value = -4995.395
$ws.Range("N132").Value = -15470.9999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H49").Value = 4500
$ws.Range("J49").Value = 4500
$ws.Range("L49").Value = 4500
$ws.Range("N49").Value = -4978

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 408.57144
$ws.Range("I16").Value = 296
$ws.Range("J16").Value = 690
$ws.Range("K16").Value = 296
$ws.Range("L16").Value = 690
$ws.Range("M16").Value = -9
$ws.Range("N16").Value = -1264
$ws.Range("H31").Value = 3995.3823
$ws.Range("I31").Value = 3495.9285
$ws.Range("J31").Value = 4345
$ws.Range("K31").Value = 3495.9285
$ws.Range("L31").Value = 4345
$ws.Range("M31").Value = -3200.9285
$ws.Range("N31").Value = -4935
$ws.Range("H34").Value = 3995.3823
$ws.Range("I34").Value = 3495.9285
$ws.Range("J34").Value = 4345
$ws.Range("K34").Value = 3495.9285
$ws.Range("L34").Value = 4345
$ws.Range("M34").Value = -3293.9285
$ws.Range("N34").Value = -4749
$ws.Range("H113").Value = 408.57144
$ws.Range("I113").Value = 296
$ws.Range("J113").Value = 690
$ws.Range("K113").Value = 296
$ws.Range("L113").Value = 690
$ws.Range("M113").Value = 1874
$ws.Range("N113").Value = -5030

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 418.18182
$ws.Range("I16").Value = 333.33334
$ws.Range("J16").Value = 800
$ws.Range("K16").Value = 1000.00002
$ws.Range("L16").Value = 2400
$ws.Range("M16").Value = -827.0000200000001
$ws.Range("N16").Value = -2746
$ws.Range("H22").Value = 145618.72
$ws.Range("I22").Value = 1888
$ws.Range("J22").Value = 169573.83
$ws.Range("K22").Value = 5664
$ws.Range("L22").Value = 508721.49
$ws.Range("M22").Value = -5495
$ws.Range("N22").Value = -509059.49
$ws.Range("H23").Value = 559.5
$ws.Range("I23").Value = 20
$ws.Range("J23").Value = 585.1905
$ws.Range("K23").Value = 60
$ws.Range("L23").Value = 1755.5715
$ws.Range("M23").Value = 175
$ws.Range("N23").Value = -2225.5715
$ws.Range("H27").Value = 145618.72
$ws.Range("I27").Value = 1888
$ws.Range("J27").Value = 169573.83
$ws.Range("K27").Value = 5664
$ws.Range("L27").Value = 508721.49
$ws.Range("M27").Value = -5562
$ws.Range("N27").Value = -508925.49
$ws.Range("H40").Value = 215.41667
$ws.Range("I40").Value = 144.09091
$ws.Range("K40").Value = 576.36364
$ws.Range("M40").Value = -507.36364
$ws.Range("H42").Value = 3020
$ws.Range("J42").Value = 3020
$ws.Range("L42").Value = 9060
$ws.Range("N42").Value = -10128
$ws.Range("H56").Value = 5150
$ws.Range("I56").Value = 5150
$ws.Range("K56").Value = 5150
$ws.Range("M56").Value = -4620
$ws.Range("H80").Value = 1000
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 1000
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 3000
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -4872
$ws.Range("H83").Value = 1000
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 1000
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 9000
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -18360
$ws.Range("H113").Value = 2551569.8
$ws.Range("I113").Value = 542.1177
$ws.Range("J113").Value = 8333899
$ws.Range("K113").Value = 1626.3531
$ws.Range("L113").Value = 25001697
$ws.Range("M113").Value = 543.6469
$ws.Range("N113").Value = -25006037
$ws.Range("H122").Value = 701.8387
$ws.Range("I122").Value = 353.375
$ws.Range("K122").Value = 3180.375
$ws.Range("M122").Value = -730.375
$ws.Range("H125").Value = 2010
$ws.Range("I125").Value = 420
$ws.Range("J125").Value = 3600
$ws.Range("K125").Value = 1260
$ws.Range("L125").Value = 10800
$ws.Range("M125").Value = 3660
$ws.Range("N125").Value = -20640

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 27.363636
$ws.Range("I2").Value = 25.166666
$ws.Range("J2").Value = 30
$ws.Range("K2").Value = 25.166666
$ws.Range("L2").Value = 30
$ws.Range("M2").Value = 87.83333400000001
$ws.Range("N2").Value = -256
$ws.Range("H70").Value = 4110.3447
$ws.Range("I70").Value = 3820.8667
$ws.Range("J70").Value = 4420.5
$ws.Range("K70").Value = 3820.8667
$ws.Range("L70").Value = 4420.5
$ws.Range("M70").Value = -3550.8667
$ws.Range("N70").Value = -4960.5
$ws.Range("H73").Value = 4110.3447
$ws.Range("I73").Value = 3820.8667
$ws.Range("J73").Value = 4420.5
$ws.Range("K73").Value = 3820.8667
$ws.Range("L73").Value = 4420.5
$ws.Range("M73").Value = -2884.8667
$ws.Range("N73").Value = -6292.5
$ws.Range("H102").Value = 2241.52
$ws.Range("I102").Value = 2641.1765
$ws.Range("J102").Value = 1392.25
$ws.Range("K102").Value = 2641.1765
$ws.Range("L102").Value = 1392.25
$ws.Range("M102").Value = -1019.1765
$ws.Range("N102").Value = -4636.25

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1574.1852
$ws.Range("I16").Value = 1670.05
$ws.Range("J16").Value = 1300.2858
$ws.Range("K16").Value = 1670.05
$ws.Range("L16").Value = 1300.2858
$ws.Range("M16").Value = -1500.05
$ws.Range("N16").Value = -1640.2858
$ws.Range("H22").Value = 696.75
$ws.Range("I22").Value = 393.5
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 393.5
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -98.5
$ws.Range("N22").Value = -1590
$ws.Range("H27").Value = 696.75
$ws.Range("I27").Value = 393.5
$ws.Range("J27").Value = 1000
$ws.Range("K27").Value = 393.5
$ws.Range("L27").Value = 1000
$ws.Range("M27").Value = -286.5
$ws.Range("N27").Value = -1214
$ws.Range("H55").Value = 292.93332
$ws.Range("I55").Value = 225.5
$ws.Range("J55").Value = 370
$ws.Range("K55").Value = 225.5
$ws.Range("L55").Value = 370
$ws.Range("M55").Value = -52.5
$ws.Range("N55").Value = -716
$ws.Range("H68").Value = 62502708
$ws.Range("I68").Value = 111112710
$ws.Range("K68").Value = 111112710
$ws.Range("M68").Value = -111111961
$ws.Range("H71").Value = 62502708
$ws.Range("I71").Value = 111112710
$ws.Range("K71").Value = 555563550
$ws.Range("M71").Value = -555559806
$ws.Range("H87").Value = 34475.6
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 34475.6
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 34475.6
$ws.Range("M87").ClearContents()
$ws.Range("N87").Value = -36721.6
$ws.Range("H90").Value = 34475.6
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 34475.6
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 103426.8
$ws.Range("M90").ClearContents()
$ws.Range("N90").Value = -114658.8

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 780
$ws.Range("I107").Value = 650
$ws.Range("J107").Value = 866.6667
$ws.Range("K107").Value = 1950
$ws.Range("L107").Value = 2600.0001
$ws.Range("M107").Value = -30
$ws.Range("N107").Value = -6440.0001
$ws.Range("H132").Value = 20837.566
$ws.Range("I132").Value = 34624.668
$ws.Range("K132").Value = 103874.004
$ws.Range("M132").Value = -101344.004
